# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.414.74'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '1.916.21'
$ws.Range("E3").Value = '  +1.07%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '0.726'
$ws.Range("E5").Value = '  +10.87%  '
$ws.Range("D6").Value = '254.16'
$ws.Range("E6").Value = '  +3.83%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").Value = '40.85'
$ws.Range("E8").Value = '  -1.57%  '
$ws.Range("D9").Value = '0.357'
$ws.Range("E9").Value = '  +1.43%  '
$ws.Range("D10").Value = '52.40'
$ws.Range("E10").Value = '  +0.50%  '
$ws.Range("D11").Value = '0.0757'
$ws.Range("E11").Value = '  +6.21%  '
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("D13").Value = '2.192.91'
$ws.Range("E13").Value = '  +1.05%  '
$ws.Range("D14").Value = '12.73'
$ws.Range("E14").Value = '  +5.75%  '
$ws.Range("D15").Value = '0.720'
$ws.Range("E15").Value = '  +3.36%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").Value = '4.94'
$ws.Range("E16").Value = '  +2.06%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '1.911.81'
$ws.Range("E17").Value = '  +1.14%  '
$ws.Range("D18").Value = '35.425.87'
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("D19").Value = '74.53'
$ws.Range("E19").Value = '  +4.40%  '
$ws.Range("D20").Value = '0.0₃0842'
$ws.Range("E20").Value = '  +2.80%  '
$ws.Range("D21").Value = '243.90'
$ws.Range("E21").Value = '  +1.50%  '
$ws.Range("D22").Value = '13.11'
$ws.Range("E22").Value = '  +5.03%  '
$ws.Range("D23").Value = '5.11'
$ws.Range("E23").Value = '  +6.51%  '
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("D25").Value = '2.49'
$ws.Range("E25").Value = '  +1.41%  '
$ws.Range("E26").Value = '  +2.47%  '
$ws.Range("D27").Value = '167.26'
$ws.Range("E27").Value = '  -1.90%  '
$ws.Range("D28").Value = '8.66'
$ws.Range("E28").Value = '  +2.19%  '
$ws.Range("D29").Value = '18.77'
$ws.Range("E29").Value = '  +2.51%  '
$ws.Range("E30").Value = '  +5.88%  '
$ws.Range("D31").Value = '4.126.21'
$ws.Range("E31").Value = '  +19.38%  '
$ws.Range("D32").Value = '4.41'
$ws.Range("E32").Value = '  +6.79%  '
$ws.Range("D33").Value = '1.99'
$ws.Range("E33").Value = '  +14.87%  '
$ws.Range("D34").Value = '1.63'
$ws.Range("E34").Value = '  +22.99%  '
$ws.Range("D35").Value = '0.0580'
$ws.Range("E35").Value = '  +2.82%  '
$ws.Range("D36").Value = '4.23'
$ws.Range("E36").Value = '  +2.89%  '
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").Value = '0.922'
$ws.Range("E38").Value = '  -1.84%  '
$ws.Range("D39").Value = '2.02'
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").Value = '17.56'
$ws.Range("E40").Value = '  +7.78%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.0218'
$ws.Range("E41").Value = '  +3.87%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '97.75'
$ws.Range("E42").Value = '  +8.80%  '
$ws.Range("E43").Value = '  +2.69%  '
$ws.Range("D44").Value = '0.0656'
$ws.Range("E44").Value = '  +1.16%  '
$ws.Range("D45").Value = '1.342.08'
$ws.Range("E45").Value = '  +0.15%  '
$ws.Range("E46").Value = '  +3.17%  '
$ws.Range("D47").Value = '2.42'
$ws.Range("E47").Value = '  +1.17%  '
$ws.Range("D48").Value = '6.75'
$ws.Range("E48").Value = '  +2.91%  '
$ws.Range("E49").Value = '  -0.21%  '
$ws.Range("D50").Value = '44.89'
$ws.Range("E50").Value = '  -5.70%  '
$ws.Range("D51").Value = '11.89'
$ws.Range("E51").Value = '  +5.98%  '
